# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns (E/F) for both rows ---
$overview.Range("E2").Value = $status
$overview.Range("F2").Value = $status
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status

# widen the status columns to fit the longer text
$overview.Range("E1").ColumnWidth = 29.9777050018311
$overview.Range("F1").ColumnWidth = 29.9777050018311

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = $status
$zhcn.Range("C3").Value = $status

$zhcn.Range("I2").Value = "71e5ab70-2d31-4b81-8b6d-0c30f10fc123.md"
$zhcn.Range("J2").Value = "71e5ab70-2d31-4b81-8b6d-0c30f10fc123.48ff3ce051c196043a2a6a21076c954f4e8a92a5.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-09 12:57:26"

$zhcn.Range("I3").Value = "c91de6a4-f413-4091-b2c1-db5858047da9.md"
$zhcn.Range("J3").Value = "c91de6a4-f413-4091-b2c1-db5858047da9.58010fc8ee31bbf9cd88a4d8a441deeea7b0a844.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-09 12:57:26"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ae2d1c0912cce7134fa5a40b825c199a1bc2286/e2e/71e5ab70-2d31-4b81-8b6d-0c30f10fc123.md", "", "", "71e5ab70-2d31-4b81-8b6d-0c30f10fc123.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ae2d1c0912cce7134fa5a40b825c199a1bc2286/e2e/c91de6a4-f413-4091-b2c1-db5858047da9.md", "", "", "c91de6a4-f413-4091-b2c1-db5858047da9.md") | Out-Null
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("I3").Font.Underline = 2
$zhcn.Range("I3").Font.Color = 15570276

$zhcn.Range("C1").ColumnWidth = 29.9777050018311
$zhcn.Range("I1").ColumnWidth = 40
$zhcn.Range("J1").ColumnWidth = 40

# --- de-de sheet ---
$dede.Range("C2").Value = $status
$dede.Range("C3").Value = $status

$dede.Range("I2").Value = "71e5ab70-2d31-4b81-8b6d-0c30f10fc123.md"
$dede.Range("J2").Value = "71e5ab70-2d31-4b81-8b6d-0c30f10fc123.48ff3ce051c196043a2a6a21076c954f4e8a92a5.de-de.xlf"
$dede.Range("K2").Value = "2016-09-09 12:57:43"

$dede.Range("I3").Value = "c91de6a4-f413-4091-b2c1-db5858047da9.md"
$dede.Range("J3").Value = "c91de6a4-f413-4091-b2c1-db5858047da9.58010fc8ee31bbf9cd88a4d8a441deeea7b0a844.de-de.xlf"
$dede.Range("K3").Value = "2016-09-09 12:57:43"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ae2d1c0912cce7134fa5a40b825c199a1bc2286/e2e/71e5ab70-2d31-4b81-8b6d-0c30f10fc123.md", "", "", "71e5ab70-2d31-4b81-8b6d-0c30f10fc123.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ae2d1c0912cce7134fa5a40b825c199a1bc2286/e2e/c91de6a4-f413-4091-b2c1-db5858047da9.md", "", "", "c91de6a4-f413-4091-b2c1-db5858047da9.md") | Out-Null
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276
$dede.Range("I3").Font.Underline = 2
$dede.Range("I3").Font.Color = 15570276

$dede.Range("C1").ColumnWidth = 29.9777050018311
$dede.Range("I1").ColumnWidth = 40
$dede.Range("J1").ColumnWidth = 40
